$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name
$ws.Name = "Through 2022-05-04"

# Update header label in B1 (shared string "May 2022 (through May 03)" -> "May 2022 (through May 04)")
$ws.Range("B1").Value = "May 2022 (through May 04)"

# Data updates (carjacking counts by neighborhood/month)
$ws.Range("G4").Value = 2        # Humboldt Park, May 2021: 1 -> 2
$ws.Range("AF5").Value = 2       # Garfield Park, May 2016: (new) -> 2
$ws.Range("B6").Value = 2        # Chicago Lawn, May 2022: 1 -> 2
$ws.Range("AA7").Value = 1       # North Lawndale, May 2017: (new) -> 1
$ws.Range("AF11").Value = 1      # Roseland, May 2016: (new) -> 1
$ws.Range("B15").Value = 1       # Lake View, May 2022: (new) -> 1
$ws.Range("AF16").Value = 1      # Little Italy, UIC, May 2016: (new) -> 1
$ws.Range("V20").Value = 1       # Woodlawn, May 2018: (new) -> 1
$ws.Range("G25").Value = 1       # Auburn Gresham, May 2021: (new) -> 1
$ws.Range("AF29").Value = 1      # West Pullman, May 2016: (new) -> 1
$ws.Range("L40").Value = 1       # Near South Side, May 2020: (new) -> 1
$ws.Range("AA48").Value = 1      # Albany Park, May 2017: (new) -> 1
$ws.Range("G89").Value = 1       # Streeterville, May 2021: (new) -> 1
